# Rename the embedded-picture "file name" labels that Word stores on the
# inline <wp:docPr>/name attribute for the logo pictures living in the
# document's headers and footers.
#
#   - BTec_Logo-Orange pictures (both headers):   image2.jpg -> image1.jpg
#   - PearsonLogo pictures      (both footers):   image1.png -> image2.png
#
# InlineShape has no settable .Name on the Word object model, but
# InlineShape.ConvertToShape() hands back a floating Shape whose .Name
# setter does write through to <wp:docPr name="...">; converting it back
# with Shape.ConvertToInlineShape() restores the original inline layout
# (<wp:inline>) so only the name actually changes.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

function Rename-LogoInlineShape($range, $newName) {
    $inline = $range.InlineShapes.Item(1)
    $shape = $inline.ConvertToShape()
    $shape.Name = $newName
    [void]$shape.ConvertToInlineShape()
}

# Headers: BTec logo, image2.jpg -> image1.jpg
for ($i = 1; $i -le 2; $i++) {
    $h = $sec.Headers.Item($i)
    if ($h.Exists -and $h.Range.InlineShapes.Count -gt 0) {
        Rename-LogoInlineShape $h.Range "image1.jpg"
    }
}

# Footers: Pearson logo, image1.png -> image2.png
for ($i = 1; $i -le 2; $i++) {
    $f = $sec.Footers.Item($i)
    if ($f.Exists -and $f.Range.InlineShapes.Count -gt 0) {
        Rename-LogoInlineShape $f.Range "image2.png"
    }
}
